$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172, shifting existing rows 172-174 down to 173-175.
$ws.Rows.Item(172).Insert()

# Populate the new row 172 with this week's price entry (same market/category,
# new date + updated volume/price figures).
$ws.Cells.Item(172, 1).Value = 8
$ws.Cells.Item(172, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 45267
$ws.Cells.Item(172, 5).Value = 4
$ws.Cells.Item(172, 6).Value = 100112028
$ws.Cells.Item(172, 7).Value = "Sandia"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 1200
$ws.Cells.Item(172, 11).Value = 700
$ws.Cells.Item(172, 12).Value = 800
$ws.Cells.Item(172, 13).Value = 750
$ws.Cells.Item(172, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(172, 15).Value = "Perú"
$ws.Cells.Item(172, 16).Value = 750
$ws.Cells.Item(172, 17).Value = 1
$ws.Cells.Item(172, 18).Value = "Hortaliza"
